# Hortaliza, Vega Central Mapocho de Santiago - Berenjena
# Weekly price-logic update: insert two new daily records (row 149 & 150)
# for date 2021-12-23, pushing the rest of the historical rows down by two
# (old row 149 -> 151, ..., old row 216 -> 218).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 149; Excel shifts everything from
# row 149 downward by one row per Insert() call, and copies the formatting
# of the row above (so column D keeps its date number format).
$ws.Rows.Item(149).Insert()
$ws.Rows.Item(149).Insert()

# New row 149: "Primera" quality entry for 2021-12-23
$ws.Cells.Item(149, 1).Value = 9
$ws.Cells.Item(149, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(149, 3).Value = "Metropolitana"
$ws.Cells.Item(149, 4).Value = 44553
$ws.Cells.Item(149, 5).Value = 13
$ws.Cells.Item(149, 6).Value = 100112001
$ws.Cells.Item(149, 7).Value = "Berenjena"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 106
$ws.Cells.Item(149, 11).Value = 7000
$ws.Cells.Item(149, 12).Value = 8000
$ws.Cells.Item(149, 13).Value = 7500
$ws.Cells.Item(149, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 150
$ws.Cells.Item(149, 17).Value = 50
$ws.Cells.Item(149, 18).Value = "Hortaliza"

# New row 150: "Segunda" quality entry for 2021-12-23
$ws.Cells.Item(150, 1).Value = 9
$ws.Cells.Item(150, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44553
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = 100112001
$ws.Cells.Item(150, 7).Value = "Berenjena"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Segunda"
$ws.Cells.Item(150, 10).Value = 61
$ws.Cells.Item(150, 11).Value = 6000
$ws.Cells.Item(150, 12).Value = 6000
$ws.Cells.Item(150, 13).Value = 6000
$ws.Cells.Item(150, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 60
$ws.Cells.Item(150, 17).Value = 100
$ws.Cells.Item(150, 18).Value = "Hortaliza"
